$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the account-number column keeps its text formatting (large numbers
# must stay as text, not be coerced to scientific-notation numbers)
$ws.Range("C2:C3").NumberFormat = "@"

# Update existing row 2
$ws.Range("A2").Value = "ALI EXPRESSE"
$ws.Range("B2").Value = "DDDD"
$ws.Range("C2").Value = "114365978523435433246652"
$ws.Range("D2").Value = "BMCE"
$ws.Range("E2").Value = "BMCE MAARIF"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/DR IFRAN"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 12000
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 10800

# Add new row 3
$ws.Range("A3").Value = "KHALID TAGHMAOUI"
$ws.Range("B3").Value = "BB132345"
$ws.Range("C3").Value = "114655862235099841255452"
$ws.Range("D3").Value = "bmce"
$ws.Range("E3").Value = "bmce"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "001/DR IFRAN"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 9000
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 16200
